$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "What is the capital of France?"
$ws.Range("A3").Value = "Who wrote One Hundred Years of Solitude?"
$ws.Range("A4").Value = "What year did man land on the Moon?"
$ws.Range("B2").Value = "Paris"
